{"js": "// Replace the 25 three-digit-by-one-digit multiplication answers in the\n// practice table with the new set of problems/answers, keeping every run's\n// formatting (font, size, paragraph alignment, etc.) untouched.\nconst replacements = [\n  [\"718\u00d75=3590\", \"919\u00d74=3676\"],\n  [\"519\u00d77=3633\", \"375\u00d73=1125\"],\n  [\"618\u00d75=3090\", \"222\u00d79=1998\"],\n  [\"137\u00d75=685\", \"794\u00d76=4764\"],\n  [\"307\u00d77=2149\", \"368\u00d72=736\"],\n  [\"676\u00d75=3380\", \"230\u00d73=690\"],\n  [\"689\u00d73=2067\", \"601\u00d76=3606\"],\n  [\"897\u00d78=7176\", \"641\u00d76=3846\"],\n  [\"157\u00d75=785\", \"561\u00d76=3366\"],\n  [\"901\u00d78=7208\", \"757\u00d78=6056\"],\n  [\"792\u00d79=7128\", \"463\u00d77=3241\"],\n  [\"809\u00d77=5663\", \"255\u00d75=1275\"],\n  [\"286\u00d77=2002\", \"603\u00d72=1206\"],\n  [\"495\u00d76=2970\", \"797\u00d76=4782\"],\n  [\"605\u00d76=3630\", \"689\u00d76=4134\"],\n  [\"858\u00d72=1716\", \"415\u00d76=2490\"],\n  [\"232\u00d77=1624\", \"257\u00d73=771\"],\n  [\"455\u00d78=3640\", \"279\u00d72=558\"],\n  [\"946\u00d78=7568\", \"889\u00d72=1778\"],\n  [\"342\u00d72=684\", \"137\u00d74=548\"],\n  [\"251\u00d76=1506\", \"943\u00d78=7544\"],\n  [\"679\u00d72=1358\", \"851\u00d78=6808\"],\n  [\"506\u00d78=4048\", \"568\u00d73=1704\"],\n  [\"684\u00d73=2052\", \"886\u00d75=4430\"],\n  [\"218\u00d73=654\", \"914\u00d75=4570\"],\n];\n\nconst body = context.document.body;\n\nfor (const [oldText, newText] of replacements) {\n  const results = body.search(oldText, { matchCase: true, matchWholeWord: false });\n  results.load(\"items\");\n  await context.sync();\n\n  if (results.items.length === 0) {\n    throw new Error(\"Text not found: \" + oldText);\n  }\n\n  for (let i = 0; i < results.items.length; i++) {\n    results.items[i].insertText(newText, Word.InsertLocation.replace);\n  }\n}\n\nawait context.sync();\n", "ps1": "# Replace the 25 three-digit-by-one-digit multiplication answers in the\n# practice table with the new set of problems/answers, keeping every run's\n# formatting (font, size, paragraph alignment, etc.) untouched.\n$d = $word.ActiveDocument\n\n$replacements = @(\n    @(\"718\u00d75=3590\", \"919\u00d74=3676\"),\n    @(\"519\u00d77=3633\", \"375\u00d73=1125\"),\n    @(\"618\u00d75=3090\", \"222\u00d79=1998\"),\n    @(\"137\u00d75=685\",  \"794\u00d76=4764\"),\n    @(\"307\u00d77=2149\", \"368\u00d72=736\"),\n    @(\"676\u00d75=3380\", \"230\u00d73=690\"),\n    @(\"689\u00d73=2067\", \"601\u00d76=3606\"),\n    @(\"897\u00d78=7176\", \"641\u00d76=3846\"),\n    @(\"157\u00d75=785\",  \"561\u00d76=3366\"),\n    @(\"901\u00d78=7208\", \"757\u00d78=6056\"),\n    @(\"792\u00d79=7128\", \"463\u00d77=3241\"),\n    @(\"809\u00d77=5663\", \"255\u00d75=1275\"),\n    @(\"286\u00d77=2002\", \"603\u00d72=1206\"),\n    @(\"495\u00d76=2970\", \"797\u00d76=4782\"),\n    @(\"605\u00d76=3630\", \"689\u00d76=4134\"),\n    @(\"858\u00d72=1716\", \"415\u00d76=2490\"),\n    @(\"232\u00d77=1624\", \"257\u00d73=771\"),\n    @(\"455\u00d78=3640\", \"279\u00d72=558\"),\n    @(\"946\u00d78=7568\", \"889\u00d72=1778\"),\n    @(\"342\u00d72=684\",  \"137\u00d74=548\"),\n    @(\"251\u00d76=1506\", \"943\u00d78=7544\"),\n    @(\"679\u00d72=1358\", \"851\u00d78=6808\"),\n    @(\"506\u00d78=4048\", \"568\u00d73=1704\"),\n    @(\"684\u00d73=2052\", \"886\u00d75=4430\"),\n    @(\"218\u00d73=654\",  \"914\u00d75=4570\")\n)\n\nforeach ($pair in $replacements) {\n    $find = $d.Content.Find\n    $find.ClearFormatting()\n    $find.Replacement.ClearFormatting()\n    $find.Text = $pair[0]\n    $find.Replacement.Text = $pair[1]\n    $find.Execute($find.Text, $false, $true, $false, $false, $false, $true, 1, $false, $find.Replacement.Text, 2) | Out-Null\n}\n"}
